$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Remove the two now-resolved bug bullets:
#  - "Le stop tourne du mauvais sens ... dans la routine;"
#  - "Si le mât s’arrête trop vite en régulation ... PWM MIN aussi."
# They are paragraphs 2 and 3; deleting paragraph 2 twice removes both,
# since each following paragraph shifts up after a delete.
$tr.Paragraphs(2, 1).Delete()
$tr.Paragraphs(2, 1).Delete()

# The "Dans les longs trajets..." bullet (now paragraph 3) used to end in
# two separate runs reading " et ça va " / "bien;" - merge them back into
# a single run (text is unchanged). Locate the merge point via the
# (ASCII, unambiguous) "msg" run that immediately precedes it.
$traj = $tr.Paragraphs(3, 1)
$msgIdx = $traj.Text.IndexOf("msg")
$mergeStart = $traj.Start + $msgIdx + "msg".Length
$etCaVaSpan = $tr.Characters($mergeStart, " et ça va bien;".Length)
$etCaVaSpan.Text = " et ça va bien;"

# The last bullet (now paragraph 4) used to end in two separate runs -
# merge the whole paragraph back into a single run (text is unchanged).
$eeprom = $tr.Paragraphs(4, 1)
$eepromSpan = $tr.Characters($eeprom.Start, $eeprom.Length)
$eepromSpan.Text = "S’il y a un power-down pendant que l’EEPROM fait une écriture interne, la valeur est erronée. Pourrait être évité avec une redondance."
